# Weekly update: insert the latest "Ajo" (garlic) price record for
# "Feria Lagunitas de Puerto Montt" above the existing history, pushing
# every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 450:465 down to 451:466, leaving a blank row 450 behind.
$ws.Rows.Item(450).Insert()

# Populate the newly inserted row 450 with this week's record.
$ws.Range("A450").Value = 4
$ws.Range("B450").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C450").Value = "Los Lagos"
$ws.Range("D450").Value = 45075
$ws.Range("E450").Value = 10
$ws.Range("F450").Value = 100112003
$ws.Range("G450").Value = "Ajo"
$ws.Range("H450").Value = "Chino"
$ws.Range("I450").Value = "Primera"
$ws.Range("J450").Value = 40
$ws.Range("K450").Value = 20000
$ws.Range("L450").Value = 20000
$ws.Range("M450").Value = 20000
$ws.Range("N450").Value = "$/caja 10 kilos"
$ws.Range("O450").Value = "China"
$ws.Range("P450").Value = 2000
$ws.Range("Q450").Value = 10
$ws.Range("R450").Value = "Hortaliza"
